# Weekly update: a new Espárragos price record for Terminal Hortofrutícola
# Agro Chillán was added. It belongs chronologically right after the
# existing row 17 (date 2022-11-18 / serial 44874), so insert a fresh row
# at position 18 — this pushes the previous rows 18-44 down to 19-45 and
# keeps the rest of the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(18).Insert()

# Populate the newly inserted row 18 with the new record's data.
$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(18, 3).Value = "Ñuble"
$ws.Cells.Item(18, 4).Value = 44902
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(18, 6).Value = 300000000
$ws.Cells.Item(18, 7).Value = "Espárragos"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 2000
$ws.Cells.Item(18, 11).Value = 900
$ws.Cells.Item(18, 12).Value = 1000
$ws.Cells.Item(18, 13).Value = 950
$ws.Cells.Item(18, 14).Value = "$/kilo"
$ws.Cells.Item(18, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(18, 16).Value = 950
$ws.Cells.Item(18, 17).Value = 1
$ws.Cells.Item(18, 18).Value = "Hortaliza"
